# Generate Report for Handoff
# Replace the old handoff-run identifiers/timestamps with the new ones
# across the Overview / zh-cn / de-de sheets, keeping each hyperlink's
# visible "display" text in sync with its cell's text (as Excel does
# when a hyperlinked cell's content is edited).

$wb = $excel.ActiveWorkbook

$oldMd  = "08ea4e39-94c5-441f-aa08-60e5f22443ac.md"
$newMd  = "c96733f2-7435-4697-aaef-c5a9950890b9.md"

$oldZh  = "08ea4e39-94c5-441f-aa08-60e5f22443ac.0127f1e37af421ae83eec80cb4dc1d502f7b606a.zh-cn.xlf"
$newZh  = "c96733f2-7435-4697-aaef-c5a9950890b9.46a260b8ab1279ac898847d4999e8ecf59615e8e.zh-cn.xlf"

$oldDe  = "08ea4e39-94c5-441f-aa08-60e5f22443ac.0127f1e37af421ae83eec80cb4dc1d502f7b606a.de-de.xlf"
$newDe  = "c96733f2-7435-4697-aaef-c5a9950890b9.46a260b8ab1279ac898847d4999e8ecf59615e8e.de-de.xlf"

$oldOverviewDate = "2016-12-14 01:12:16"
$newOverviewDate = "2016-13-14 01:13:21"

$oldZhDate = "2016-03-14 01:09:35"
$newZhDate = "2016-03-14 01:13:17"

$oldDeDate = "2016-03-14 01:12:16"
$newDeDate = "2016-03-14 01:13:21"

function Update-HyperlinkDisplay($ws, $oldText, $newText) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.TextToDisplay() -eq $oldText) {
            $hl.TextToDisplay = $newText
        }
    }
}

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMd
$wsOverview.Range("D2").Value = $newOverviewDate
Update-HyperlinkDisplay $wsOverview $oldMd $newMd

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newMd
$wsZh.Range("D2").Value = $newZh
$wsZh.Range("E2").Value = $newZhDate
Update-HyperlinkDisplay $wsZh $oldMd $newMd
Update-HyperlinkDisplay $wsZh $oldZh $newZh

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newMd
$wsDe.Range("D2").Value = $newDe
$wsDe.Range("E2").Value = $newDeDate
Update-HyperlinkDisplay $wsDe $oldMd $newMd
Update-HyperlinkDisplay $wsDe $oldDe $newDe
